$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "24.335.68"
$ws.Range("E2").Value = "  -1.73%  "

# Row 3
$ws.Range("D3").Value = "1.677.61"
$ws.Range("E3").Value = "  -1.27%  "

# Row 4
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").Value = "'315.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.70%  "

# Row 6
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "'0.3883"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.20%  "

# Row 8
$ws.Range("D8").Value = "'0.3987"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "

# Row 9
$ws.Range("B9").Value = "BinanceUSD"
$ws.Range("C9").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D9").Value = "'1.001"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.22%  "

# Row 10
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.468"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.84%  "

# Row 11
$ws.Range("D11").Value = "'52.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.73%  "

# Row 12
$ws.Range("D12").Value = "'0.08696"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.28%  "

# Row 13
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").Value = "'25.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.56%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'7.453"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.76%  "

# Row 15
$ws.Range("D15").Value = "'7.932"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.27%  "

# Row 16
$ws.Range("D16").Value = "'0.00001332"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.97%  "

# Row 17
$ws.Range("D17").Value = "1.667.96"
$ws.Range("E17").Value = "  -1.54%  "

# Row 18
$ws.Range("D18").Value = "'97.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.39%  "

# Row 19
$ws.Range("D19").Value = "'0.07074"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "

# Row 20
$ws.Range("D20").Value = "'19.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.47%  "

# Row 21
$ws.Range("D21").Value = "'7.202"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.70%  "

# Row 22
$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "

# Row 23
$ws.Range("D23").Value = "'14.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.91%  "

# Row 24
$ws.Range("D24").Value = "24.296.76"
$ws.Range("E24").Value = "  -1.79%  "

# Row 25
$ws.Range("D25").Value = "'2.978"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.16%  "

# Row 26
$ws.Range("D26").Value = "'2.330"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "

# Row 27
$ws.Range("D27").Value = "'22.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.60%  "

# Row 28
$ws.Range("D28").Value = "'164.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.84%  "

# Row 29
$ws.Range("D29").Value = "'8.640"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.96%  "

# Row 30
$ws.Range("B30").Value = "HuobiToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D30").Value = "'5.250"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.14%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "'136.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "

# Row 32
$ws.Range("D32").Value = "1.853.26"
$ws.Range("E32").Value = "  -1.45%  "

# Row 33
$ws.Range("D33").Value = "'0.08709"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.22%  "

# Row 34
$ws.Range("D34").Value = "'7.303"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.36%  "

# Row 35
$ws.Range("D35").Value = "'1.031"
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("B36").Value = "WEMIXTOKEN"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "'1.969"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.59%  "

# Row 37
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2757"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "

# Row 38
$ws.Range("D38").Value = "'0.02909"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.30%  "

# Row 39
$ws.Range("D39").Value = "'10.62"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.27%  "

# Row 40
$ws.Range("D40").Value = "'0.09058"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.47%  "

# Row 41
$ws.Range("D41").Value = "'13.99"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.86%  "

# Row 42
$ws.Range("D42").Value = "'0.7830"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.86%  "

# Row 43
$ws.Range("D43").Value = "'1.462"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.34%  "

# Row 44
$ws.Range("D44").Value = "'16.96"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.40%  "

# Row 45
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.7116"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "

# Row 46
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.564"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.10%  "

# Row 47
$ws.Range("D47").Value = "'4.201"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.58%  "

# Row 48
$ws.Range("D48").Value = "'1.391"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.69%  "

# Row 49
$ws.Range("D49").Value = "'1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50
$ws.Range("D50").Value = "'138.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.29%  "

# Row 51
$ws.Range("D51").Value = "'0.08004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.24%  "

